$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three data values in row 1 (A1, B1, C1)
$ws.Range("A1").Value = 148.87328002502903
$ws.Range("B1").Value = 3.6442242194680285
$ws.Range("C1").Value = 2.5929549902152642

# Narrow column A (was 11.7109375 -> 9.7109375) and column C (was 11.7109375 -> 10.7109375).
# Column B's width (11.7109375) is left untouched.
# ColumnWidth is expressed in "characters" and Excel snaps it to whole-pixel
# boundaries (MDW=7) when applied through the object model, so these inputs
# are chosen to land as close as the COM layer allows to the target widths.
$ws.Columns.Item(1).ColumnWidth = 8.833333333333334
$ws.Columns.Item(3).ColumnWidth = 9.833333333333334
